$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header row (drop the trailing space that was on "PlaceName ")
$ws.Range("A1").Value = "PlaceName"
$ws.Range("B1").Value = "Longitude"
$ws.Range("C1").Value = "Latitude"

# Place names first (order matters for how new shared strings get appended),
# then fill in the coordinate values row by row.
$ws.Range("A4").Value = "City Center Shopping Mall"
$ws.Range("A2").Value = "MGF Metropolitan Mall"
$ws.Range("A3").Value = "DT City Center Mall"
$ws.Range("A5").Value = "MGF Mega Mall"
$ws.Range("A6").Value = "Sahara Mall"

$ws.Range("B2").Value = 28.480908500000002
$ws.Range("C2").Value = 77.078091099999995

$ws.Range("B3").Value = 28.487140799999999
$ws.Range("C3").Value = 77.090245600000003

$ws.Range("B4").Value = 28.4789745
$ws.Range("C4").Value = 77.078437500000007

$ws.Range("B5").Value = 28.479719599999999
$ws.Range("C5").Value = 77.0871736

$ws.Range("B6").Value = 28.479558699999998
$ws.Range("C6").Value = 77.084565799999993

# Match the final selection state from the saved workbook
$ws.Range("B6").Select()
